$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11

$t.Cell(1,1).Range.Text = "92 x 31${nl}  3    1${nl}  ----${nl}9|    |${nl}2|    |"
$t.Cell(1,2).Range.Text = "12 x 18${nl}  1    8${nl}  ----${nl}1|    |${nl}2|    |"
$t.Cell(1,3).Range.Text = "18 x 48${nl}  4    8${nl}  ----${nl}1|    |${nl}8|    |"

$t.Cell(2,1).Range.Text = "86 x 61${nl}  6    1${nl}  ----${nl}8|    |${nl}6|    |"
$t.Cell(2,2).Range.Text = "17 x 70${nl}  7    0${nl}  ----${nl}1|    |${nl}7|    |"
$t.Cell(2,3).Range.Text = "59 x 33${nl}  3    3${nl}  ----${nl}5|    |${nl}9|    |"

$t.Cell(3,1).Range.Text = "38 x 92${nl}  9    2${nl}  ----${nl}3|    |${nl}8|    |"
$t.Cell(3,2).Range.Text = "49 x 99${nl}  9    9${nl}  ----${nl}4|    |${nl}9|    |"
$t.Cell(3,3).Range.Text = "87 x 19${nl}  1    9${nl}  ----${nl}8|    |${nl}7|    |"

$t.Cell(4,1).Range.Text = "93 x 38${nl}  3    8${nl}  ----${nl}9|    |${nl}3|    |"
$t.Cell(4,2).Range.Text = "19 x 18${nl}  1    8${nl}  ----${nl}1|    |${nl}9|    |"
$t.Cell(4,3).Range.Text = "27 x 60${nl}  6    0${nl}  ----${nl}2|    |${nl}7|    |"

$t.Cell(5,1).Range.Text = "71 x 69${nl}  6    9${nl}  ----${nl}7|    |${nl}1|    |"
$t.Cell(5,2).Range.Text = "14 x 93${nl}  9    3${nl}  ----${nl}1|    |${nl}4|    |"
$t.Cell(5,3).Range.Text = "21 x 72${nl}  7    2${nl}  ----${nl}2|    |${nl}1|    |"
